$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that are no longer part of the dataset (rows 19-30)
$ws.Range("A19:C30").EntireRow.Delete()

# Update headers
$ws.Range("B1").Value = "X"
$ws.Range("C1").Value = "Y"

# Update data rows 2..18 (column A keeps its existing 0..16 values)
for ($i = 0; $i -le 16; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $i
    $ws.Cells.Item($row, 3).Value = (5 * $i) + 3
}
